$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core scenario inputs change (pool sizes 1000/1000/996 -> 500/500/506) ---
$ws.Range("C9").Value = 500
$ws.Range("C10").Value = 500
$ws.Range("C11").Value = 506

# --- I26 now branches depending on whether CEX price (C11) is above or below
#     the V4 pool price (C9), instead of always assuming C11 > C9 ---
$ws.Range("I26").Formula = "=IF(C11>C9,((I6-F8)/(F7-F8))*(1-(2*C6)),((F7-I6)/(F7-F8))*(1-(2*C6)))"

# --- New footnote explaining #NUM/DIV0 errors near the bottom of the model.
#     B19/B20 are written before D14 so the shared-string table picks up the
#     same ordering as the authored workbook (long sentence first, "*" last). ---
$ws.Range("B19").Value = "* If there's a #NUM or DIV/0 error it means CEX/DEX "
$ws.Range("B20").Value = "arbitrage has no expectation of profit."
$ws.Range("D14").Value = "*"

# Footnote text uses a small italic font (size 10)
$ws.Range("B19").Font.Italic = $true
$ws.Range("B19").Font.Size = 10
$ws.Range("B20").Font.Italic = $true
$ws.Range("B20").Font.Size = 10

# --- Column L widened slightly to fit the updated numbers ---
$ws.Columns("L").ColumnWidth = 18.17

# --- Selection cursor moved from D17 to G14 ---
$ws.Range("G14").Select()
